$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card14")

# --- Row 15: fill "nan" into previously blank cells B15:K15, M15 ---
$ws.Range("B15:K15").Value = 'nan'
$ws.Range("M15").Value = 'nan'

# --- Row 16: new event row ---
# A16 = "14" (kept as text, matching column A convention in this sheet)
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = '14'

# B16:K16 and M16 stay blank but present (match style of an existing blank-style data cell)
$ws.Range("B16:K16").Style = $ws.Range("B2").Style
$ws.Range("M16").Style = $ws.Range("B2").Style

# L16 = date-like text (backslash-separated, not Excel-date-parseable, stays text naturally)
$ws.Range("L16").Value = '5\1\2024'

# N16, O16 = Arabic text values
$ws.Range("N16").Value = 'تم سن الفلاتس +تغيير جريده 1 + سن السليندر'
$ws.Range("O16").Value = 'الخبير'
